# Slide 4 ("Diffuse Lighting & Shadows"):
#   1. Speaker notes gain a new bullet about lighting entities.
#   2. The "tutorial" caption is reworded from "Based on" to "Expanded on".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(4)

# --- 1. Notes: append "Apply lighting to entities as well." as a new paragraph ---
# NotesPage.Shapes.AddPlaceholder(2) returns the existing notes-body
# placeholder (ppPlaceholderBody) rather than creating a duplicate one.
$notesBody = $slide.NotesPage.Shapes.AddPlaceholder(2)
$notesRange = $notesBody.TextFrame.TextRange
$notesRange.Text = $notesRange.Text + "`r" + "Apply lighting to entities as well."

# --- 2. Body placeholder: "Based on the tutorial we followed" -> "Expanded on ..." ---
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Based on the tutorial we followed") {
            $shape.TextFrame.TextRange.Text = "Expanded on the tutorial we followed"
        }
    }
}
